$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wrong SMILES code for "Frumin" (row 3) - it was a duplicate/incorrect
# value; replace it with the correct SMILES string.
$ws.Range("B3").Value = "CCOP(=S)(OCC)SCCSCC"

# Add a hyperlink on the corrected SMILES cell pointing to a similarity search.
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.chemeo.com/similar?smiles=CCOP%28%3dS%29%28OCC%29SCCSCC", [Type]::Missing, [Type]::Missing, "https://www.chemeo.com/similar?smiles=CCOP%28%3dS%29%28OCC%29SCCSCC")

# The rows that used the "code" font style (s=2) for plain, non-SMILES-looking
# values now get the regular/default text style (same as column A) - copy the
# formatting (not the value) from A2 onto B2, B3, B4, B5, B6. Done after the
# hyperlink is added so the auto-applied "Hyperlink" style gets overwritten
# back to the plain style, matching the source file (which keeps B3 unstyled).
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)

# Match the saved selection position recorded in the workbook.
$ws.Range("F9").Select()
